# Update relevant transaction data on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# rate
$ws.Range("E2").Value = 1.2624
# usdValue
$ws.Range("F2").Value = 126.24
# accountNumber
$ws.Range("K2").Value = 1234956578
